# Add the 2020 column (N) to the indicator sheet, mirroring the existing
# 2010-2019 columns (D:M) that already hold one column per year.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing year column (M) into the new
# column (N) so the new cells pick up the same styles as the other year
# cells (header style for row 4, data style for row 5), then overwrite
# with the new year's values.
$ws.Range("M4").Copy($ws.Range("N4"))
$ws.Range("N4").Value = 2020

$ws.Range("M5").Copy($ws.Range("N5"))
$ws.Range("N5").Value = 2.1

# Leave the selection where the author's saved view had it.
$ws.Range("N9").Select()
